$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1103.3334
$ws.Range("I9").Value = 155
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 155
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 14
$ws.Range("N9").Value = -3338

$ws.Range("H21").Value = 9500
$ws.Range("I21").Value = 6000
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -5532
$ws.Range("N21").Value = -10936

$ws.Range("H23").Value = 9500
$ws.Range("I23").Value = 6000
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -5766
$ws.Range("N23").Value = -10468

$ws.Range("H29").Value = 838
$ws.Range("I29").Value = 95
$ws.Range("J29").Value = 1333.3334
$ws.Range("K29").Value = 285
$ws.Range("L29").Value = 4000.0002
$ws.Range("M29").Value = -4
$ws.Range("N29").Value = -4562.0002

$ws.Range("H38").Value = 361.2
$ws.Range("I38").Value = 201.5
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 604.5
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -232.5
$ws.Range("N38").Value = -3744

$ws.Range("H43").Value = 1228.4546
$ws.Range("I43").Value = 1755.4615
$ws.Range("J43").Value = 1065.3334
$ws.Range("K43").Value = 1755.4615
$ws.Range("L43").Value = 1065.3334
$ws.Range("M43").Value = -1686.4615
$ws.Range("N43").Value = -1203.3334

$ws.Range("H58").Value = 628.0769
$ws.Range("I58").Value = 96.25
$ws.Range("J58").Value = 864.44446
$ws.Range("K58").Value = 288.75
$ws.Range("L58").Value = 2593.33338
$ws.Range("M58").Value = -138.75
$ws.Range("N58").Value = -2893.33338

$ws.Range("H113").Value = 3114.6296
$ws.Range("I113").Value = 4098.75
$ws.Range("J113").Value = 2943.4783
$ws.Range("K113").Value = 4098.75
$ws.Range("L113").Value = 2943.4783
$ws.Range("M113").Value = -844.75
$ws.Range("N113").Value = -9451.478300000001

$ws.Range("H116").Value = 2256.111
$ws.Range("I116").Value = 2001
$ws.Range("J116").Value = 2575
$ws.Range("K116").Value = 2001
$ws.Range("L116").Value = 2575
$ws.Range("M116").Value = 1441
$ws.Range("N116").Value = -9459

$ws.Range("H132").Value = 10210380
$ws.Range("I132").Value = 16135594
$ws.Range("J132").Value = 5846.4443
$ws.Range("K132").Value = 48406782
$ws.Range("L132").Value = 17539.3329
$ws.Range("M132").Value = -48404252
$ws.Range("N132").Value = -22599.3329

$ws.Range("H133").Value = 39980
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39980
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39980
$ws.Range("N133").Value = -50100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 787.3099999999999
$ws.Range("I32").Value = 775.7474
$ws.Range("J32").Value = 1007
$ws.Range("K32").Value = 775.7474
$ws.Range("L32").Value = 1007
$ws.Range("M32").Value = -488.7474

$ws.Range("H61").Value = 1272.439
$ws.Range("I61").Value = 949.0625
$ws.Range("J61").Value = 2422.2222
$ws.Range("K61").Value = 949.0625
$ws.Range("L61").Value = 2422.2222
$ws.Range("M61").Value = -737.0625
$ws.Range("N61").Value = -2846.2222

$ws.Range("H132").Value = 1732226.9
$ws.Range("I132").Value = 1756.75
$ws.Range("J132").Value = 9807754
$ws.Range("K132").Value = 5270.25
$ws.Range("L132").Value = 29423262
$ws.Range("M132").Value = -2740.25
$ws.Range("N132").Value = -29428322

$ws.Range("H136").Value = 1272.439
$ws.Range("I136").Value = 949.0625
$ws.Range("J136").Value = 2422.2222
$ws.Range("K136").Value = 2847.1875
$ws.Range("L136").Value = 7266.6666
$ws.Range("M136").Value = -297.1875
$ws.Range("N136").Value = -12366.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 500
$ws.Range("I16").Value = 500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -330
$ws.Range("N16").ClearContents()

$ws.Range("H42").Value = 180000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 180000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 180000
$ws.Range("N42").Value = -180656

$ws.Range("H99").Value = 1606.4651
$ws.Range("I99").Value = 1336.7142
$ws.Range("J99").Value = 2110
$ws.Range("K99").Value = 1336.7142
$ws.Range("L99").Value = 2110
$ws.Range("M99").Value = 161.2858000000001
$ws.Range("N99").Value = -5106

$ws.Range("H134").Value = 3006805.2
$ws.Range("I134").Value = 1073.4073
$ws.Range("J134").Value = 11122281
$ws.Range("K134").Value = 3220.2219
$ws.Range("L134").Value = 33366843
$ws.Range("M134").Value = -685.2219000000005
$ws.Range("N134").Value = -33371913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 41673904
$ws.Range("I99").Value = 66675756
$ws.Range("J99").Value = 4144.4443
$ws.Range("K99").Value = 66675756
$ws.Range("L99").Value = 4144.4443
$ws.Range("M99").Value = -66674258
$ws.Range("N99").Value = -7140.4443

$ws.Range("H105").Value = 3595.5845
$ws.Range("I105").Value = 3583.2876
$ws.Range("J105").Value = 3820
$ws.Range("K105").Value = 3583.2876
$ws.Range("L105").Value = 3820
$ws.Range("M105").Value = -1836.2876
$ws.Range("N105").Value = -7314

$ws.Range("H126").Value = 41673904
$ws.Range("I126").Value = 66675756
$ws.Range("J126").Value = 4144.4443
$ws.Range("K126").Value = 200027268
$ws.Range("L126").Value = 12433.3329
$ws.Range("M126").Value = -200024798
$ws.Range("N126").Value = -17373.3329

$ws.Range("H132").Value = 7248141
$ws.Range("I132").Value = 1364.0834
$ws.Range("J132").Value = 15153716
$ws.Range("K132").Value = 4092.2502
$ws.Range("L132").Value = 45461148
$ws.Range("M132").Value = -1562.2502
$ws.Range("N132").Value = -45466208

$ws.Range("H134").Value = 23810886
$ws.Range("I134").Value = 1474.4375
$ws.Range("J134").Value = 100001000
$ws.Range("K134").Value = 4423.3125
$ws.Range("L134").Value = 300003000
$ws.Range("M134").Value = -1888.3125
$ws.Range("N134").Value = -300008070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6017.278
$ws.Range("I68").Value = 428.2857
$ws.Range("J68").Value = 9573.909
$ws.Range("K68").Value = 1284.8571
$ws.Range("L68").Value = 28721.727
$ws.Range("M68").Value = -473.8571000000002
$ws.Range("N68").Value = -30343.727

$ws.Range("H71").Value = 6017.278
$ws.Range("I71").Value = 428.2857
$ws.Range("J71").Value = 9573.909
$ws.Range("K71").Value = 3854.5713
$ws.Range("L71").Value = 86165.181
$ws.Range("M71").Value = 201.4286999999999
$ws.Range("N71").Value = -94277.181

$ws.Range("H131").Value = 781.09
$ws.Range("I131").Value = 404
$ws.Range("J131").Value = 832.51135
$ws.Range("K131").Value = 1212
$ws.Range("L131").Value = 2497.53405
$ws.Range("M131").Value = 3828
$ws.Range("N131").Value = -12577.53405

$ws.Range("H139").Value = 185012.38
$ws.Range("I139").Value = 1066.0769
$ws.Range("J139").Value = 334468.75
$ws.Range("K139").Value = 3198.2307
$ws.Range("L139").Value = 1003406.25
$ws.Range("M139").Value = 1941.7693
$ws.Range("N139").Value = -1013686.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 12000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 12000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12336

$ws.Range("H132").Value = 10999.435
$ws.Range("I132").Value = 9315
$ws.Range("J132").Value = 12837
$ws.Range("K132").Value = 27945
$ws.Range("L132").Value = 38511
$ws.Range("M132").Value = -25415
$ws.Range("N132").Value = -43571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18187486
$ws.Range("I132").Value = 38463430
$ws.Range("J132").Value = 9052
$ws.Range("K132").Value = 115390290
$ws.Range("L132").Value = 27156
$ws.Range("M132").Value = -115387760
$ws.Range("N132").Value = -32216

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H118").Value = 20392
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 20392
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 20392
$ws.Range("N118").Value = -23706

$ws.Range("H122").Value = 14655.436
$ws.Range("I122").Value = 22186.584
$ws.Range("J122").Value = 2605.6
$ws.Range("K122").Value = 66559.75199999999
$ws.Range("L122").Value = 7816.799999999999
$ws.Range("M122").Value = -64109.75199999999

$ws.Range("H132").Value = 18404.016
$ws.Range("I132").Value = 20449.303
$ws.Range("J132").Value = 10065.538
$ws.Range("K132").Value = 61347.909
$ws.Range("L132").Value = 30196.614
$ws.Range("M132").Value = -58817.909
$ws.Range("N132").Value = -35256.614

$ws.Range("H139").Value = 115430
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 115430
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 115430
$ws.Range("N139").Value = -125710
$ws.Range("M139").ClearContents()
